$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 / row 24 are the highlighted "cut-line" rows. Their last visible column
# (E, since F was blank) did not need a right-hand border. Now that column F
# (Competition) is getting real content, give E7/E24 the same right-edge
# border that the (soon to be replaced) F7 formatting already has, so the
# highlighted row is still properly boxed in.
$ws.Range("F7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null

# Give the whole Competition column (F) the same formatting already used by
# F2, then fill in the competition name for each of the two leaderboard
# blocks.
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F3:F35").PasteSpecial(-4122) | Out-Null

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = "Single Stableford"
}
for ($r = 19; $r -le 35; $r++) {
    $ws.Cells.Item($r, 6).Value = "Test Stabelford"
}
